$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -11.5308
$ws.Range("C7").Value = -11.8086
$ws.Range("B8").Value = 4.763100000000001
$ws.Range("A12").Value = -22.79460000000001
$ws.Range("B12").Value = 5.617800000000001
$ws.Range("B14").Value = 9.561800000000007
$ws.Range("C19").Value = -13.57429999999999
$ws.Range("D19").Value = -8.283099999999996
$ws.Range("C21").Value = -12.9953
$ws.Range("B22").Value = 4.799300000000005
$ws.Range("C24").Value = -11.4175
